$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(448).Insert()

$ws.Cells.Item(448, 1).Value = 5
$ws.Cells.Item(448, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(448, 3).Value = "Maule"
$ws.Cells.Item(448, 4).Value = 45166
$ws.Cells.Item(448, 5).Value = 7
$ws.Cells.Item(448, 6).Value = 100112008
$ws.Cells.Item(448, 7).Value = "Coliflor"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Primera"
$ws.Cells.Item(448, 10).Value = 3000
$ws.Cells.Item(448, 11).Value = 800
$ws.Cells.Item(448, 12).Value = 800
$ws.Cells.Item(448, 13).Value = 800
$ws.Cells.Item(448, 14).Value = "`$/unidad"
$ws.Cells.Item(448, 15).Value = "Región del Maule"
$ws.Cells.Item(448, 16).Value = 800
$ws.Cells.Item(448, 17).Value = 1
$ws.Cells.Item(448, 18).Value = "Hortaliza"
